# Label on top of practice trials
# Uppercase the respModal values in column C (rows 2-25), widen column C
# to match column B, and move the selection to C26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    $cell.Value = $val.ToUpper()
}

$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

$ws.Range("C26").Select()
